$wb = $excel.ActiveWorkbook

# --- 1. Update selection / view on the existing "amt_923_tk" sheet ---
$ws5 = $wb.Worksheets.Item("amt_923_tk")
$ws5.Select()
$ws5.Range("E1:G1").Select()

# --- 2. Add the new "amt_929_tk" sheet after the last sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws6 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws6.Name = "amt_929_tk"

# --- 3. Header row ---
$ws6.Range("A1").Value = 'entryCode'
$ws6.Range("B1").Value = 'total'
$ws6.Range("C1").Value = 'total.valid'
$ws6.Range("D1").Value = 'Worker.ID'
$ws6.Range("E1").Value = 'toRate'
$ws6.Range("F1").Value = 'Comment'
$ws6.Range("G1").Value = 'turker.Index'

# --- 4. Data rows, written bottom-up (row 36 -> row 2), column D before column A, ---
#        to reproduce the exact shared-string insertion order of the original commit.
# row 36
$ws6.Range("D36").Value = 'A3CXJ18H14DXBN'
$ws6.Range("A36").Value = '866945fe96a91e48d993ec72c6038ed3'
$ws6.Range("B36").Value = 8
$ws6.Range("C36").Value = 2
$ws6.Range("E36").Value = 0

# row 35
$ws6.Range("D35").Value = 'A2NMETDZGRDQ55'
$ws6.Range("A35").Value = '133e1eedfe840b9dfcb610e14858a136'
$ws6.Range("B35").Value = 8
$ws6.Range("C35").Value = 2
$ws6.Range("E35").Value = 0

# row 34
$ws6.Range("D34").Value = 'A1C0H8G0YI15MN'
$ws6.Range("A34").Value = 'd3b1ec38087526478e14d653b20ab1e9'
$ws6.Range("B34").Value = 8
$ws6.Range("C34").Value = 3
$ws6.Range("E34").Value = 0

# row 33
$ws6.Range("D33").Value = 'A2NB54BDTLA0QS'
$ws6.Range("A33").Value = '036da1d114403dcd4911a56723081951'
$ws6.Range("B33").Value = 8
$ws6.Range("C33").Value = 3
$ws6.Range("E33").Value = 0

# row 32
$ws6.Range("D32").Value = 'A3MBLDIFREDC5M'
$ws6.Range("A32").Value = '72c9dfb85f0e65a593a9fc3b4dc11435'
$ws6.Range("B32").Value = 8
$ws6.Range("C32").Value = 4
$ws6.Range("E32").Value = 0

# row 31
$ws6.Range("D31").Value = 'A2H74STYH3PLT1'
$ws6.Range("A31").Value = '36de6316630cf172bde72c9242216ed4'
$ws6.Range("B31").Value = 8
$ws6.Range("C31").Value = 7
$ws6.Range("E31").Value = 0

# row 30
$ws6.Range("D30").Value = 'ADKSME3J8B0PA'
$ws6.Range("A30").Value = 'ff1f954c7711f01be600490f94dfa1f6'
$ws6.Range("B30").Value = 8
$ws6.Range("C30").Value = 8
$ws6.Range("E30").Value = 1
$ws6.Range("G30").Value = 183

# row 29
$ws6.Range("D29").Value = 'AMA426ARRZO1O'
$ws6.Range("A29").Value = 'e75be82b2a5ccf056c6654819bc05c87'
$ws6.Range("B29").Value = 8
$ws6.Range("C29").Value = 8
$ws6.Range("E29").Value = 1
$ws6.Range("G29").Value = 182

# row 28
$ws6.Range("D28").Value = 'A2DLH5XGBNYXWS'
$ws6.Range("A28").Value = 'e33b1cd62152a988de5c0257c29fdc2d'
$ws6.Range("B28").Value = 8
$ws6.Range("C28").Value = 8
$ws6.Range("E28").Value = 1
$ws6.Range("G28").Value = 181

# row 27
$ws6.Range("D27").Value = 'A1IMYG6LYQHOPD'
$ws6.Range("A27").Value = 'd2e71d2029d85d61d6ff06f16b91c6c1'
$ws6.Range("B27").Value = 8
$ws6.Range("C27").Value = 8
$ws6.Range("E27").Value = 1
$ws6.Range("G27").Value = 180

# row 26
$ws6.Range("D26").Value = 'A1CE2XPYCDRHVZ'
$ws6.Range("A26").Value = 'cc3d483ba59e741337b189f96ec96d4b'
$ws6.Range("B26").Value = 8
$ws6.Range("C26").Value = 8
$ws6.Range("E26").Value = 1
$ws6.Range("G26").Value = 179

# row 25
$ws6.Range("D25").Value = 'A3F832E9XKFRRP'
$ws6.Range("A25").Value = 'b96dece8ee1f8645b6482675b6c653ae'
$ws6.Range("B25").Value = 8
$ws6.Range("C25").Value = 8
$ws6.Range("E25").Value = 1
$ws6.Range("G25").Value = 178

# row 24
$ws6.Range("D24").Value = 'A16U1RK5OHN08F'
$ws6.Range("A24").Value = 'b5b6cabacae3055ee30422cdbaaac221'
$ws6.Range("B24").Value = 8
$ws6.Range("C24").Value = 8
$ws6.Range("E24").Value = 1
$ws6.Range("G24").Value = 177

# row 23
$ws6.Range("D23").Value = 'A1HBIE5LRTQK1L'
$ws6.Range("A23").Value = 'b0319f326020de32ad3af45b08c55930'
$ws6.Range("B23").Value = 8
$ws6.Range("C23").Value = 8
$ws6.Range("E23").Value = 1
$ws6.Range("G23").Value = 176

# row 22
$ws6.Range("D22").Value = 'A220217I0IXYX3'
$ws6.Range("A22").Value = 'a6d4ad452fde42bcaff7286a0477024f'
$ws6.Range("B22").Value = 8
$ws6.Range("C22").Value = 8
$ws6.Range("E22").Value = 1
$ws6.Range("G22").Value = 175

# row 21
$ws6.Range("D21").Value = 'A1OBB3PWYWK9KK'
$ws6.Range("A21").Value = '9e8886135ca79ca16b880f944fa69acc'
$ws6.Range("B21").Value = 8
$ws6.Range("C21").Value = 8
$ws6.Range("E21").Value = 1
$ws6.Range("G21").Value = 174
$ws6.Range("A21").NumberFormat = "0.00E+00"

# row 20
$ws6.Range("D20").Value = 'A3UQRXQB7BQ9FE'
$ws6.Range("A20").Value = '8956e6c21c8cc34b0745eb4f3151ce4f'
$ws6.Range("B20").Value = 8
$ws6.Range("C20").Value = 8
$ws6.Range("E20").Value = 1
$ws6.Range("G20").Value = 173

# row 19
$ws6.Range("D19").Value = 'A3QXV94C2J0LQ'
$ws6.Range("A19").Value = '88fd91371fcf76200f68a615f4a359de'
$ws6.Range("B19").Value = 8
$ws6.Range("C19").Value = 8
$ws6.Range("E19").Value = 1
$ws6.Range("G19").Value = 172

# row 18
$ws6.Range("D18").Value = 'A2GZ0MWS800M6R'
$ws6.Range("A18").Value = '798785112a4b2872e4ae00fb6b1eec2f'
$ws6.Range("B18").Value = 8
$ws6.Range("C18").Value = 8
$ws6.Range("E18").Value = 0
$ws6.Range("E18").Font.Color = 255

# row 17
$ws6.Range("D17").Value = 'AV08UM669CO02'
$ws6.Range("A17").Value = '63ee40e6f5823c5076e408eafec13701'
$ws6.Range("B17").Value = 8
$ws6.Range("C17").Value = 8
$ws6.Range("E17").Value = 1
$ws6.Range("G17").Value = 171

# row 16
$ws6.Range("D16").Value = 'A1G4ZQ8NAB94TT'
$ws6.Range("A16").Value = '625878bbe7175a773b52b1e66b1d2105'
$ws6.Range("B16").Value = 8
$ws6.Range("C16").Value = 8
$ws6.Range("E16").Value = 1
$ws6.Range("G16").Value = 170

# row 15
$ws6.Range("D15").Value = 'AOLLFPCWXJVA6'
$ws6.Range("A15").Value = '61137aa97a7daefe8c0b5b7a6e7e66c7'
$ws6.Range("B15").Value = 8
$ws6.Range("C15").Value = 8
$ws6.Range("E15").Value = 0
$ws6.Range("E15").Font.Color = 255

# row 14
$ws6.Range("D14").Value = 'A1DUH3RLI00YQM'
$ws6.Range("A14").Value = '5f697dbe33cc16a5676ae1956d0f86aa'
$ws6.Range("B14").Value = 8
$ws6.Range("C14").Value = 8
$ws6.Range("E14").Value = 1
$ws6.Range("G14").Value = 169

# row 13
$ws6.Range("D13").Value = 'A2ZLDAQZIN5WZC'
$ws6.Range("A13").Value = '5f59e63618bb40a726c83a1da9fabcf4'
$ws6.Range("B13").Value = 8
$ws6.Range("C13").Value = 8
$ws6.Range("E13").Value = 1
$ws6.Range("G13").Value = 168

# row 12
$ws6.Range("D12").Value = 'A2C4271VBZQCR7'
$ws6.Range("A12").Value = '5bbd81263d7ac05504e0dd0e3279b79c'
$ws6.Range("B12").Value = 8
$ws6.Range("C12").Value = 8
$ws6.Range("E12").Value = 1
$ws6.Range("G12").Value = 167

# row 11
$ws6.Range("D11").Value = 'A1FHBPFI2UNASX'
$ws6.Range("A11").Value = '4a47c45e5d5f36706ebec7a4f163c9a8'
$ws6.Range("B11").Value = 8
$ws6.Range("C11").Value = 8
$ws6.Range("E11").Value = 1
$ws6.Range("G11").Value = 166

# row 10
$ws6.Range("D10").Value = 'A1XB03X4J35ATE'
$ws6.Range("A10").Value = '49bf7822e28fc4244ff4b480647bbc27'
$ws6.Range("B10").Value = 8
$ws6.Range("C10").Value = 8
$ws6.Range("E10").Value = 1
$ws6.Range("G10").Value = 165

# row 9
$ws6.Range("D9").Value = 'A39NKZDUFD70NV'
$ws6.Range("A9").Value = '462343b842422e8203c07aa7c24b86a9'
$ws6.Range("B9").Value = 8
$ws6.Range("C9").Value = 8
$ws6.Range("E9").Value = 1
$ws6.Range("G9").Value = 164

# row 8
$ws6.Range("D8").Value = 'A3LUXDAIWTYTZL'
$ws6.Range("A8").Value = '433cee823dc7f7e921a5096abab454b0'
$ws6.Range("B8").Value = 8
$ws6.Range("C8").Value = 8
$ws6.Range("E8").Value = 0
$ws6.Range("E8").Font.Color = 255

# row 7
$ws6.Range("D7").Value = 'A2809ZZ59YT4C0'
$ws6.Range("A7").Value = '3ddc21c75fa5950646a46a99416ab3a7'
$ws6.Range("B7").Value = 8
$ws6.Range("C7").Value = 8
$ws6.Range("E7").Value = 1
$ws6.Range("G7").Value = 163

# row 6
$ws6.Range("D6").Value = 'A1QDWIDKBZK759'
$ws6.Range("A6").Value = '3cca3f67f9aff39473637c4ca1b3b8f4'
$ws6.Range("B6").Value = 8
$ws6.Range("C6").Value = 8
$ws6.Range("E6").Value = 1
$ws6.Range("G6").Value = 162

# row 5
$ws6.Range("D5").Value = 'AK7LGB1QOGA1P'
$ws6.Range("A5").Value = '342e7519b04c94956408b39a8300e1c2'
$ws6.Range("B5").Value = 8
$ws6.Range("C5").Value = 8
$ws6.Range("E5").Value = 1
$ws6.Range("G5").Value = 161
$ws6.Range("A5").NumberFormat = "0.00E+00"

# row 4
$ws6.Range("D4").Value = 'A4ZW4GNQ98HV6'
$ws6.Range("A4").Value = '3129db589e16da52605c903eac08e2af'
$ws6.Range("B4").Value = 8
$ws6.Range("C4").Value = 8
$ws6.Range("E4").Value = 1
$ws6.Range("G4").Value = 160

# row 3
$ws6.Range("D3").Value = 'A2PE32I58ANCDD'
$ws6.Range("A3").Value = '20a55ea9098263a79f5eacbd2e93e3c9'
$ws6.Range("B3").Value = 8
$ws6.Range("C3").Value = 8
$ws6.Range("E3").Value = 1
$ws6.Range("G3").Value = 159

# row 2
$ws6.Range("D2").Value = 'A2J9S68Y0ROJ8W'
$ws6.Range("A2").Value = '1808b557f62e5112da80067517465799'
$ws6.Range("B2").Value = 8
$ws6.Range("C2").Value = 8
$ws6.Range("E2").Value = 1
$ws6.Range("G2").Value = 158

# --- 5. Comment text (F column) for the invalid rows; written last to match string order ---
$ws6.Range("F8").Value = 'video playing is choppy'
$ws6.Range("F15").Value = 'video playing is choppy'
$ws6.Range("F18").Value = 'camera has not been started'

# --- 6. Final selection & active sheet ---
$ws6.Range("E18").Select()
$ws6.Select()

Write-Output "done"